# Workbook / sheet handles
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert three new rows right above the current row 7 (CathayPacific),
#    shifting CathayPacific..United from rows 7-23 down to rows 10-26.
# ---------------------------------------------------------------------------
$ws.Range("A7:A9").EntireRow.Insert()

# ---------------------------------------------------------------------------
# 2. Populate the three newly inserted rows (Asiana, British, Cargolux).
# ---------------------------------------------------------------------------
$ws.Range("A7").Value = "Asiana"
$ws.Range("B7").Value = "988"
$ws.Range("C7").Value = "in progress"
$ws.Range("D7").Value = "https://www.asianacargo.com/tracking/viewTraceAirWaybill.do?lang=en"

$ws.Range("A8").Value = "British"
$ws.Range("B8").Value = "125"
$ws.Range("C8").Value = "in progress"
$ws.Range("D8").Value = "https://www.iagcargo.com/en/home"

$ws.Range("A9").Value = "Cargolux"
$ws.Range("B9").Value = "172"
$ws.Range("C9").Value = "in progress"
$ws.Range("D9").Value = "https://cvtnt.champ.aero/trackntrace"

# ---------------------------------------------------------------------------
# 3. Append a new row 27 (VirginAtlantic) after United (now row 26).
# ---------------------------------------------------------------------------
$ws.Range("A27").Value = "VirginAtlantic"
$ws.Range("B27").Value = "932"
$ws.Range("C27").Value = "in progress"
$ws.Range("D27").Value = "https://cargo.virgin-atlantic.com/gb/en/track/track-your-cargo.html?prefix=932&number=56409673&track=go"

# ---------------------------------------------------------------------------
# 4. Rebuild every hyperlink. Inserting rows does not renumber the
#    worksheet's <hyperlink> anchors automatically, so clear them all and
#    re-add them in the exact order the target workbook uses (this also
#    drives the r:id numbering of the relationships part).
# ---------------------------------------------------------------------------
$ws.Cells.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("D2"),  "https://www.airbridgecargo.com/en/tracking/")
$ws.Hyperlinks.Add($ws.Range("D3"),  "https://www.airtahitinui.com/us-en/online-cargo-tracking")
$ws.Hyperlinks.Add($ws.Range("D4"),  "https://www.afklcargo.com/WW/en/local/app/index.jsp", "/tntsinglesearch")
$ws.Hyperlinks.Add($ws.Range("D10"), "http://www.cathaypacificcargo.com/ManageYourShipment/TrackYourShipment/tabid/108/SingleAWBNo/160-05480334-/language/en-US/Default.aspx")
$ws.Hyperlinks.Add($ws.Range("D11"), "https://cargo.china-airlines.com/CCNetv2/content/manage/ShipmentTracking.aspx?")
$ws.Hyperlinks.Add($ws.Range("D12"), "https://aviationcargo.dhl.com/aviationcargo/track/")
$ws.Hyperlinks.Add($ws.Range("D13"), "https://skychain.emirates.com/skychain/app?service=page/nwp:Trackshipmt&amp;initial=y")
$ws.Hyperlinks.Add($ws.Range("D14"), "http://www.brcargo.com/ec_web/Default.aspx?Parm2=191&amp;Parm3=undefined")
$ws.Hyperlinks.Add($ws.Range("D15"), "http://www.jal.co.jp/en/jalcargo/inter/awb/")
$ws.Hyperlinks.Add($ws.Range("D17"), "https://www.afklcargo.com/WW/en/local/app/index.jsp", "/tntsinglesearch")
$ws.Hyperlinks.Add($ws.Range("D16"), "https://lufthansa-cargo.com/eservices/etracking")
$ws.Hyperlinks.Add($ws.Range("D18"), "https://cargo.koreanair.com/en/tracking?")
$ws.Hyperlinks.Add($ws.Range("D19"), "http://www.maskargo.com/online_awb_info/index.php")
$ws.Hyperlinks.Add($ws.Range("D21"), "https://freight.qantas.com/online-tracking.html?")
$ws.Hyperlinks.Add($ws.Range("D22"), "http://www.qrcargo.com/trackshipment")
$ws.Hyperlinks.Add($ws.Range("D23"), "http://www.siacargo.com/ccn/ShipmentTrack.aspx")
$ws.Hyperlinks.Add($ws.Range("D24"), "https://www.skyteam.com/en/cargo/track-shipment/")
$ws.Hyperlinks.Add($ws.Range("D26"), "https://www.unitedcargo.com/OurNetwork/TrackingCargo1512/Tracking.jsp")
$ws.Hyperlinks.Add($ws.Range("D25"), "https://www.turkishcargo.com.tr/en/online-services/shipment-tracking")
$ws.Hyperlinks.Add($ws.Range("D5"),  "https://www.aacargo.com/AACargo/tracking")
$ws.Hyperlinks.Add($ws.Range("D6"),  "https://mycargo.amerijet.com/tracking")
$ws.Hyperlinks.Add($ws.Range("D20"), "https://www.anacargo.jp/en/int/")
$ws.Hyperlinks.Add($ws.Range("D7"),  "https://www.asianacargo.com/tracking/viewTraceAirWaybill.do?lang=en")
$ws.Hyperlinks.Add($ws.Range("D8"),  "https://www.iagcargo.com/en/home")
$ws.Hyperlinks.Add($ws.Range("D27"), "https://cargo.virgin-atlantic.com/gb/en/track/track-your-cargo.html?prefix=932&number=56409673&track=go")
$ws.Hyperlinks.Add($ws.Range("D9"),  "https://cvtnt.champ.aero/trackntrace")

# ---------------------------------------------------------------------------
# 5. Re-apply the "Hyperlink" cell style to every linked D cell (Hyperlinks.Add
#    otherwise leaves the cell with a plain default style), then restore
#    vertical-center alignment on the subset of rows that use it.
# ---------------------------------------------------------------------------
$hyperlinkCells = "D2","D3","D4","D5","D6","D7","D8","D9","D10","D11","D12","D13","D14","D15","D16","D17","D18","D19","D20","D21","D22","D23","D24","D25","D26","D27"
foreach ($addr in $hyperlinkCells) {
    $ws.Range($addr).Style = "Hyperlink"
}

$centeredCells = "D3","D10","D18","D19","D21","D23","D24","D25","D26"
foreach ($addr in $centeredCells) {
    $ws.Range($addr).VerticalAlignment = -4108
}

# ---------------------------------------------------------------------------
# 6. Selection, matching the workbook's saved cursor position.
# ---------------------------------------------------------------------------
$ws.Range("D9").Select()
